$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

# Header row
$ws2.Range("A1").Value = "P20 Roll No."
$ws2.Range("B1").Value = "Name"
$ws2.Range("C1").Value = "Total"
$ws2.Range("D1").Value = "Percentage"

# Header formatting (column A only): bold, size 12, centered + middle, wrap
$hdrA = $ws2.Range("A1")
$hdrA.Font.Bold = $true
$hdrA.Font.Size = 12
$hdrA.HorizontalAlignment = -4108
$hdrA.VerticalAlignment = -4108
$hdrA.WrapText = $true

# Data rows
$ws2.Range("A2").Value = "P20001"
$ws2.Range("B2").Value = "Student1 "
$ws2.Range("C2").Formula = "=VLOOKUP(A2,Sheet1!A1:L20,8,0)"
$ws2.Range("D2").Formula = "=VLOOKUP(A2,Sheet1!A2:L20,9,0)"
$ws2.Range("A3").Value = "P20002"
$ws2.Range("B3").Value = "Student2"
$ws2.Range("C3").Formula = "=VLOOKUP(A3,Sheet1!A2:L21,8,0)"
$ws2.Range("D3").Formula = "=VLOOKUP(A3,Sheet1!A3:L21,9,0)"
$ws2.Range("A4").Value = "P20003"
$ws2.Range("B4").Value = "Student3"
$ws2.Range("C4").Formula = "=VLOOKUP(A4,Sheet1!A3:L22,8,0)"
$ws2.Range("D4").Formula = "=VLOOKUP(A4,Sheet1!A4:L22,9,0)"
$ws2.Range("A5").Value = "P20004"
$ws2.Range("B5").Value = "Student4"
$ws2.Range("C5").Formula = "=VLOOKUP(A5,Sheet1!A4:L23,8,0)"
$ws2.Range("D5").Formula = "=VLOOKUP(A5,Sheet1!A5:L23,9,0)"
$ws2.Range("A6").Value = "P20005"
$ws2.Range("B6").Value = "Student5"
$ws2.Range("C6").Formula = "=VLOOKUP(A6,Sheet1!A5:L24,8,0)"
$ws2.Range("D6").Formula = "=VLOOKUP(A6,Sheet1!A6:L24,9,0)"
$ws2.Range("A7").Value = "P20006"
$ws2.Range("B7").Value = "Student6"
$ws2.Range("C7").Formula = "=VLOOKUP(A7,Sheet1!A6:L25,8,0)"
$ws2.Range("D7").Formula = "=VLOOKUP(A7,Sheet1!A7:L25,9,0)"
$ws2.Range("A8").Value = "P20007"
$ws2.Range("B8").Value = "Student7"
$ws2.Range("C8").Formula = "=VLOOKUP(A8,Sheet1!A7:L26,8,0)"
$ws2.Range("D8").Formula = "=VLOOKUP(A8,Sheet1!A8:L26,9,0)"
$ws2.Range("A9").Value = "P20008"
$ws2.Range("B9").Value = "Student8"
$ws2.Range("C9").Formula = "=VLOOKUP(A9,Sheet1!A8:L27,8,0)"
$ws2.Range("D9").Formula = "=VLOOKUP(A9,Sheet1!A9:L27,9,0)"
$ws2.Range("A10").Value = "P20009"
$ws2.Range("B10").Value = "Student9"
$ws2.Range("C10").Formula = "=VLOOKUP(A10,Sheet1!A9:L28,8,0)"
$ws2.Range("D10").Formula = "=VLOOKUP(A10,Sheet1!A10:L28,9,0)"
$ws2.Range("A11").Value = "P20010"
$ws2.Range("B11").Value = "Student10"
$ws2.Range("C11").Formula = "=VLOOKUP(A11,Sheet1!A10:L29,8,0)"
$ws2.Range("D11").Formula = "=VLOOKUP(A11,Sheet1!A11:L29,9,0)"
$ws2.Range("A12").Value = "P20011"
$ws2.Range("B12").Value = "Student11"
$ws2.Range("C12").Formula = "=VLOOKUP(A12,Sheet1!A11:L30,8,0)"
$ws2.Range("D12").Formula = "=VLOOKUP(A12,Sheet1!A12:L30,9,0)"
$ws2.Range("A13").Value = "P20012"
$ws2.Range("B13").Value = "Student12"
$ws2.Range("C13").Formula = "=VLOOKUP(A13,Sheet1!A12:L31,8,0)"
$ws2.Range("D13").Formula = "=VLOOKUP(A13,Sheet1!A13:L31,9,0)"
$ws2.Range("A14").Value = "P20013"
$ws2.Range("B14").Value = "Student13"
$ws2.Range("C14").Formula = "=VLOOKUP(A14,Sheet1!A13:L32,8,0)"
$ws2.Range("D14").Formula = "=VLOOKUP(A14,Sheet1!A14:L32,9,0)"
$ws2.Range("A15").Value = "P20014"
$ws2.Range("B15").Value = "Student14"
$ws2.Range("C15").Formula = "=VLOOKUP(A15,Sheet1!A14:L33,8,0)"
$ws2.Range("D15").Formula = "=VLOOKUP(A15,Sheet1!A15:L33,9,0)"
$ws2.Range("A16").Value = "P20015"
$ws2.Range("B16").Value = "Student15"
$ws2.Range("C16").Formula = "=VLOOKUP(A16,Sheet1!A15:L34,8,0)"
$ws2.Range("D16").Formula = "=VLOOKUP(A16,Sheet1!A16:L34,9,0)"
$ws2.Range("A17").Value = "P20016"
$ws2.Range("B17").Value = "Student16"
$ws2.Range("C17").Formula = "=VLOOKUP(A17,Sheet1!A16:L35,8,0)"
$ws2.Range("D17").Formula = "=VLOOKUP(A17,Sheet1!A17:L35,9,0)"
$ws2.Range("A18").Value = "P20017"
$ws2.Range("B18").Value = "Student17"
$ws2.Range("C18").Formula = "=VLOOKUP(A18,Sheet1!A17:L36,8,0)"
$ws2.Range("D18").Formula = "=VLOOKUP(A18,Sheet1!A18:L36,9,0)"
$ws2.Range("A19").Value = "P20018"
$ws2.Range("B19").Value = "Student18"
$ws2.Range("C19").Formula = "=VLOOKUP(A19,Sheet1!A18:L37,8,0)"
$ws2.Range("D19").Formula = "=VLOOKUP(A19,Sheet1!A19:L37,9,0)"
$ws2.Range("A20").Value = "P20020"
$ws2.Range("B20").Value = "Student19"
$ws2.Range("C20").Formula = "=VLOOKUP(A20,Sheet1!A19:L38,8,0)"
$ws2.Range("D20").Formula = "=VLOOKUP(A20,Sheet1!A20:L38,9,0)"

# Data rows formatting (column A only): size 12, middle vertical align, wrap (no bold, no horizontal center)
$dataA = $ws2.Range("A2:A20")
$dataA.Font.Size = 12
$dataA.VerticalAlignment = -4108
$dataA.WrapText = $true

# Column widths
$ws2.Columns.Item(2).ColumnWidth = 13
$ws2.Columns.Item(3).ColumnWidth = 21.3
$ws2.Columns.Item(4).ColumnWidth = 27.8

# Row heights
$ws2.Rows.Item(1).RowHeight = 31.5
$ws2.Rows.Item("2:20").RowHeight = 15.75

# Page setup (to emit pageSetup element)
$ws2.PageSetup.Orientation = 1

# Selection / view
$ws2.Range("H9").Select()

# Sheet1 selection should move to L20, and Sheet1 should no longer be the tabSelected sheet
$ws1.Range("L20").Select()

# Activate Sheet2 last so it is the active/selected tab
$ws2.Activate()

Write-Output "done"
